# Update symbol list (price + volume-label refresh), per
# "Updated symbol list on Thu Dec 22 07:43:28 UTC 2022 with GitHub Actions"
#
# Source cells in column D store prices as plain text (not numbers), so a
# leading apostrophe is used when assigning to force Excel to keep them as
# text instead of auto-converting to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'245.85"
$ws.Range("D3").Value  = "'22.86"
$ws.Range("D4").Value  = "'5.270"
$ws.Range("D5").Value  = "'0.05732"
$ws.Range("D6").Value  = "'3.444"
$ws.Range("D7").Value  = "'0.8098"
$ws.Range("D8").Value  = "'0.8773"
$ws.Range("D10").Value = "'0.07372"
$ws.Range("D11").Value = "'0.03017"
$ws.Range("D12").Value = "'0.03105"
$ws.Range("D13").Value = "'0.09392"
$ws.Range("D15").Value = "'0.04807"
$ws.Range("D16").Value = "'0.0005842"
$ws.Range("D18").Value = "'0.005096"
$ws.Range("D19").Value = "'0.0009968"
$ws.Range("D21").Value = "'3.750"
$ws.Range("D22").Value = "'6.295"
$ws.Range("D23").Value = "'2.191"
$ws.Range("D26").Value = "'4.160"
$ws.Range("D27").Value = "'0.0003002"
$ws.Range("D40").Value = "'0.03899"

$ws.Range("D41").Value = "'0.006774"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").Value = "'0.1066"
$ws.Range("D43").Value = "'0.003199"
$ws.Range("D44").Value = "'0.007381"
$ws.Range("D45").Value = "'0.00005639"

$ws.Range("D47").Value = "'0.6002"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D48").Value = "'0.1735"
$ws.Range("D49").Value = "'0.00002099"
